# "Generate Report for handoff"
# Updates the per-language handoff-status sheets (zh-cn, de-de) with the
# freshly generated handoff info: status text, new "Latest Handoff File"
# hyperlink + its datetime, and the handoff reason.

$wb = $excel.ActiveWorkbook

$commit = "f920b2958f06df9313379dcbfb7b0847305773ee"
$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/$commit"
$uuid = "c5c317c1-99b9-4fdc-9fbb-b71ed502eb84"
$tag = "e05eda47f1f1d6240c059956c0bbc44c516f8730"

function Update-LangSheet {
    param(
        [string]$SheetName,
        [string]$XlfFile,
        [string]$HandoffDateTime
    )

    $ws = $wb.Worksheets.Item($SheetName)

    # Re-create the existing hyperlinks plus the new one, in left-to-right /
    # row order, so the new "Latest Handoff File" link lands between the
    # two that were already there.
    $ws.Hyperlinks.Delete()

    $ws.Range("B2").Value = "Not yet handed off"

    $ws.Range("C2").Value = $XlfFile
    $ws.Range("D2").Value = $HandoffDateTime
    $ws.Range("H2").Value = "Include"

    $ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/e2e/$uuid.md", "", "", "$uuid.md")
    $ws.Hyperlinks.Add($ws.Range("C2"), "$repoBase/e2e/Loc/$XlfFile", "", "", $XlfFile)
    $ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/.localization-config", "", "", ".localization-config")
}

Update-LangSheet "zh-cn" "$uuid.$tag.zh-cn.xlf" "2016-01-08 10:39:45"
Update-LangSheet "de-de" "$uuid.$tag.de-de.xlf" "2016-01-08 10:39:54"

# The "Status" text is a shared string, so the Overview roll-up (which
# points at the very same cell text) picks up the new wording too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Not yet handed off"
$wsOverview.Range("C2").Value = "Not yet handed off"
